$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the erroneous value in the input/output table -----------------
# The last cell of row 7 ("G7") held a wrong pattern; correct it in place
# before the row is shifted below (note the intentional trailing space,
# matching the authored fix).
$ws.Range("G7").Value = "0 1 0 1 0 1 "

# --- Capture the current (pre-shift) contents of row 7 ------------------
$b7 = $ws.Range("B7").Value()
$c7 = $ws.Range("C7").Value()
$d7 = $ws.Range("D7").Value()
$e7 = $ws.Range("E7").Value()
$f7 = $ws.Range("F7").Value()
$g7 = $ws.Range("G7").Value()

# --- Give the new trailing blank cell (H7) the same look as the other ---
# --- blank / bordered cells in the table (e.g. G4) before we write it ---
$ws.Range("G4").Copy()
$ws.Range("H7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Shift the row-7 values one column to the right ----------------------
$ws.Range("H7").Value = $null
$ws.Range("G7").Value = $f7
$ws.Range("F7").Value = $e7
$ws.Range("E7").Value = $d7
$ws.Range("D7").Value = $c7
$ws.Range("C7").Value = $b7

# --- The corrected value now becomes the new first cell of the row ------
# --- (it carries no direct formatting, unlike its neighbours) -----------
$ws.Range("B7").Value = $g7
$ws.Range("B7").ClearFormats()

# --- Move the active selection -------------------------------------------
$ws.Range("C14").Select() | Out-Null
